# Schedule.xlsx — "Card Display" sheet update
#  - Update the header date string
#  - Remove the duplicated J:R block (columns J through R)
#  - Add the position / employee roster in column A, rows 7-25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card Display")

# 1. Update the title text (was "Saturday March 5, 2016")
$ws.Range("A2").Value = "Friday March 4, 2016"

# 2. Remove the duplicate card block that lived in columns J:R
$ws.Columns("J:R").Delete()

# 3. Populate the roster list in column A
$ws.Range("A7").Value  = "Bin Filler"
$ws.Range("A8").Value  = " Doris Reynolds"
$ws.Range("A9").Value  = " Zabada Mohammed"

$ws.Range("A11").Value = "Forklift"
$ws.Range("A12").Value = " George Dunn"
$ws.Range("A13").Value = " Don Coles"
$ws.Range("A14").Value = " George C Brown"

$ws.Range("A16").Value = "Line Operator"
$ws.Range("A17").Value = " Parveen Gopal"

$ws.Range("A19").Value = "QC"
$ws.Range("A20").Value = " Isabel Roseen"

$ws.Range("A22").Value = "Non Rotational"
$ws.Range("A23").Value = " Elaine Roseen"
$ws.Range("A24").Value = " Janeanne Reiswig"
$ws.Range("A25").Value = " Sandra Martin"

$ws.Range("A5").Select() | Out-Null
